# Logged Week 15 and simulated Week 16
# Appends the new week's per-play logs to the running shared-string play
# logs on YDS / ST, and updates the season-to-date aggregate totals on
# OFF / DEF / ST / TURNS / PEN to reflect the two additional games.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# YDS — append new per-play yardage logs (Week 15 + Week 16)
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("YDS")

$ws.Range("B2").Value2 = $ws.Range("B2").Value2 + " 3 -1 7 8 4 13 17 2 3 2 0 2 3 10 7 9 12 2 5 0 0 2 1 -2 1 12 3 1"
$ws.Range("C2").Value2 = $ws.Range("C2").Value2 + " 1 6 5 4 5 5 28 4 5 0 5 1 -3 0 6 6 3 5 1 3 2 0"
$ws.Range("B3").Value2 = $ws.Range("B3").Value2 + " 9 11 4 2 19 4 9 9 6 5 11 3 3 6 7 15 9 3 7 3 10 16 10 1 11 7 9 7"
$ws.Range("C3").Value2 = $ws.Range("C3").Value2 + " 8 2 16 4 15 0 3 4 9 15 9 7 7 8 14 7 7 10 36"

# ---------------------------------------------------------------------
# OFF — season-to-date offensive totals
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("OFF")

$ws.Range("C2").Value2 = 190
$ws.Range("D2").Value2 = 10
$ws.Range("E2").Value2 = 3
$ws.Range("F2").Value2 = 54
$ws.Range("G2").Value2 = 47
$ws.Range("I2").Value2 = 5
$ws.Range("J2").Value2 = 26
$ws.Range("N2").Value2 = 19

$ws.Range("C3").Value2 = 159
$ws.Range("E3").Value2 = 32
$ws.Range("F3").Value2 = 117
$ws.Range("G3").Value2 = 37
$ws.Range("H3").Value2 = 18
$ws.Range("I3").Value2 = 50
$ws.Range("J3").Value2 = 58
$ws.Range("L3").Value2 = 338
$ws.Range("M3").Value2 = 232
$ws.Range("Q3").Value2 = 567

# ---------------------------------------------------------------------
# DEF — season-to-date defensive totals
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("DEF")

$ws.Range("C2").Value2 = 160
$ws.Range("F2").Value2 = 54
$ws.Range("H2").Value2 = 5
$ws.Range("J2").Value2 = 20
$ws.Range("N2").Value2 = 16
$ws.Range("O2").Value2 = 18
$ws.Range("P2").Value2 = 10

$ws.Range("B3").Value2 = 16
$ws.Range("C3").Value2 = 168
$ws.Range("E3").Value2 = 34
$ws.Range("F3").Value2 = 93
$ws.Range("G3").Value2 = 31
$ws.Range("H3").Value2 = 26
$ws.Range("I3").Value2 = 52
$ws.Range("J3").Value2 = 50
$ws.Range("L3").Value2 = 295
$ws.Range("M3").Value2 = 177
$ws.Range("Q3").Value2 = 503

# ---------------------------------------------------------------------
# ST — special-teams totals + per-game logs
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("ST")

$ws.Range("B2").Value2 = 84
$ws.Range("D2").Value2 = 54
$ws.Range("F2").Value2 = 342
$ws.Range("J2").Value2 = 170
$ws.Range("K2").Value2 = 161
$ws.Range("L2").Value2 = 80
$ws.Range("M2").Value2 = 64

$ws.Range("B3").Value2 = 70

$ws.Range("B4").Value2 = $ws.Range("B4").Value2 + " 56"
$ws.Range("B5").Value2 = $ws.Range("B5").Value2 + " 15"
$ws.Range("B6").Value2 = $ws.Range("B6").Value2 + " 21 29"
$ws.Range("D3").Value2 = $ws.Range("D3").Value2 + " 32 36 50 42"
$ws.Range("D4").Value2 = $ws.Range("D4").Value2 + " 0 0 0 0"
$ws.Range("D5").Value2 = $ws.Range("D5").Value2 + " 0 0"

# ---------------------------------------------------------------------
# TURNS — season-to-date turnover totals
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("TURNS")

$ws.Range("C3").Value2 = 13
$ws.Range("D3").Value2 = 8
$ws.Range("E3").Value2 = 4

# ---------------------------------------------------------------------
# PEN — season-to-date penalty totals
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("PEN")

$ws.Range("B2").Value2 = 16
$ws.Range("D2").Value2 = 11
$ws.Range("D4").Value2 = 9
